$d = $word.ActiveDocument

function Split-At([string]$needle) {
    # Re-splits the run containing $needle from its neighbours by toggling
    # Bold on then off again on that exact sub-range; this engine merges
    # same-formatted adjacent runs on every text edit, so nudging (and
    # reverting) a character property is what forces a fresh run boundary
    # at both ends of the matched text.
    $r = $d.Content
    $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Bold = 1
    $r.Bold = 0
}

# ---------------------------------------------------------------------
# 1. "Set_likelihood_and_prior" + ": add section " -> merge into one run
#    (the following "for country " run must stay a separate, untouched
#    run, so re-split it off again afterwards)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Set_likelihood_and_prior: add section ", $true, $false, $false, $false, $false, $true, 1, $false, "Set_likelihood_and_prior: add section ", 2)
Split-At "for country "

# ---------------------------------------------------------------------
# 2. "Save_data_manager" (proofErr-wrapped) + ": " -> merge into one run
#    (also drops the now-unneeded proofErr spell-check markers). This is
#    the last run in its paragraph, so no re-split is needed afterwards.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Save_data_manager: ", $true, $false, $false, $false, $false, $true, 1, $false, "Save_data_manager: ", 2)

# ---------------------------------------------------------------------
# 3. "Right now, using South Africa's engagement data and Kenya's
#    disengagement data" -> "For some countries," / " using South
#    Africa's engagement data and Kenya's disengagement data". This run
#    is alone in its paragraph, so one split suffices.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Right now, using South Africa’s engagement data and Kenya’s disengagement data", $true, $false, $false, $false, $false, $true, 1, $false, "For some countries, using South Africa’s engagement data and Kenya’s disengagement data", 2)
Split-At " using South Africa’s engagement data and Kenya’s disengagement data"

# ---------------------------------------------------------------------
# 4. "if using something other than Kenya's data, update function calls
#    and add files to country-specific folders " -> three runs: "if
#    using something other than Kenya" / " or SA" / "'s data, update
#    function calls and add files to country-specific folders " (the
#    preceding "Get_suppression_rebound_data: " run must stay separate
#    and untouched)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("if using something other than Kenya’s data, update function calls and add files to country-specific folders ", $true, $false, $false, $false, $false, $true, 1, $false, "if using something other than Kenya or SA’s data, update function calls and add files to country-specific folders ", 2)
Split-At "Get_suppression_rebound_data: "
Split-At " or SA’s data, update function calls and add files to country-specific folders "
Split-At " or SA"

# ---------------------------------------------------------------------
# 5. "Right now, using Kenya's" + " suppression/" + "unsuppression"
#    (proofErr-wrapped) + " data" -> four runs: "For some countries" /
#    ", using Kenya" / " or SA's" / " suppression/unsuppression data"
#    (the trailing "; uses a function call..." run must stay separate
#    and untouched)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Right now, using Kenya’s suppression/unsuppression data", $true, $false, $false, $false, $false, $true, 1, $false, "For some countries, using Kenya or SA’s suppression/unsuppression data", 2)
Split-At "; uses a function call where all other locations assume Kenya – nothing to update "
Split-At " suppression/unsuppression data"
Split-At ", using Kenya"

# ---------------------------------------------------------------------
# 6. Remove the whole "Source new prior " bullet paragraph (the list
#    item right after "Source code:"), including its paragraph mark, so
#    "Source code:" is immediately followed by "Load new data manager".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Source new prior ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p = $rng.Paragraphs.First
$pRange = $p.Range
[void]$pRange.MoveEnd(1, 1)
$delRange = $d.Range($pRange.Start, $pRange.End)
[void]$delRange.Delete()
